# Remove the manually-typed "section number + tab" prefix runs from every
# heading paragraph (Heading 1-4). This corresponds to setting
# number_sections: FALSE in the bookdown/csasdown YAML, which stops the
# renderer from prepending a literal section number (e.g. "1.2.1") and a
# tab character before the heading text.

$d = $word.ActiveDocument

# First pass: find the paragraphs that need editing (by index) without
# mutating the document, so indices/offsets stay stable while scanning.
$targets = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -like "Heading*") {
        $t = $p.Range.Text
        if ($t -match "^[0-9]+(\.[0-9]+)*\t") {
            [void]$targets.Add($i)
        }
    }
}

# Second pass: actually delete the "number + tab" prefix from each
# matched heading paragraph.
foreach ($i in $targets) {
    $p = $d.Paragraphs($i)
    $pStart = $p.Range.Start
    $count = $p.Range.Characters.Count
    $tabIndex = -1
    for ($j = 1; $j -le $count; $j++) {
        $ch = $p.Range.Characters($j)
        if ($ch.Text -eq "`t") {
            $tabIndex = $j
            break
        }
    }
    if ($tabIndex -ge 0) {
        $tabChar = $p.Range.Characters($tabIndex)
        $absEnd = $pStart + $tabChar.End
        $delRange = $d.Range($pStart, $absEnd)
        $delRange.Delete()
    }
}
